# Auto-generated edit script: update market-price derived columns (H-N)
# across multiple sheets per the scheduled-runner price refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 996.45
$ws.Range("J28").Value = 1107.5714
$ws.Range("L28").Value = 1107.5714
$ws.Range("N28").Value = -2077.5714
$ws.Range("H34").Value = 4999.75
$ws.Range("I34").Value = 4999.75
$ws.Range("K34").Value = 4999.75
$ws.Range("M34").Value = -4796.75
$ws.Range("H36").Value = 4999.75
$ws.Range("I36").Value = 4999.75
$ws.Range("K36").Value = 4999.75
$ws.Range("M36").Value = -4284.75
$ws.Range("H62").Value = 7549.2144
$ws.Range("I62").Value = 6852
$ws.Range("K62").Value = 6852
$ws.Range("M62").Value = -6228
$ws.Range("H65").Value = 7549.2144
$ws.Range("I65").Value = 6852
$ws.Range("K65").Value = 34260
$ws.Range("M65").Value = -31140
$ws.Range("H76").Value = 4300
$ws.Range("I76").Value = 4300
$ws.Range("K76").Value = 4300
$ws.Range("M76").Value = -3985
$ws.Range("H79").Value = 4300
$ws.Range("I79").Value = 4300
$ws.Range("K79").Value = 4300
$ws.Range("M79").Value = -3208
$ws.Range("H100").Value = 4611.769
$ws.Range("I100").Value = 3746.0833
$ws.Range("K100").Value = 3746.0833
$ws.Range("M100").Value = -3205.0833
$ws.Range("H106").Value = 1825.4783
$ws.Range("I106").Value = 1557.6666
$ws.Range("K106").Value = 1557.6666
$ws.Range("M106").Value = -926.6666
$ws.Range("H107").Value = 1320.85
$ws.Range("I107").Value = 1456.1666
$ws.Range("K107").Value = 1456.1666
$ws.Range("M107").Value = 463.8334
$ws.Range("H111").Value = 1308.7368
$ws.Range("I111").Value = 1136.3077
$ws.Range("K111").Value = 3408.9231
$ws.Range("M111").Value = -341.9231
$ws.Range("H112").Value = 13482.8
$ws.Range("J112").Value = 16603.625
$ws.Range("L112").Value = 49810.875
$ws.Range("N112").Value = -52026.875
$ws.Range("H113").Value = 5195.5835
$ws.Range("I113").Value = 4053
$ws.Range("K113").Value = 4053
$ws.Range("M113").Value = -799
$ws.Range("H116").Value = 11863.1
$ws.Range("I116").Value = 16977.3
$ws.Range("J116").Value = 6748.9
$ws.Range("K116").Value = 16977.3
$ws.Range("L116").Value = 6748.9
$ws.Range("M116").Value = -13535.3
$ws.Range("N116").Value = -13632.9
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("H138").Value = 2053503.4
$ws.Range("J138").Value = 3545131
$ws.Range("L138").Value = 10635393
$ws.Range("N138").Value = -10645673

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23416.775
$ws.Range("I32").Value = 24042.51
$ws.Range("K32").Value = 24042.51
$ws.Range("M32").Value = -23755.51
$ws.Range("H43").Value = 18646.445
$ws.Range("J43").Value = 18434.5
$ws.Range("L43").Value = 18434.5
$ws.Range("N43").Value = -19060.5
$ws.Range("H112").Value = 29428.285
$ws.Range("J112").Value = 29428.285
$ws.Range("L112").Value = 29428.285
$ws.Range("N112").Value = -32382.285
$ws.Range("H139").Value = 123239.8
$ws.Range("J139").Value = 123239.8
$ws.Range("L139").Value = 123239.8
$ws.Range("N139").Value = -133519.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 21203476
$ws.Range("I11").Value = 31800214
$ws.Range("J11").Value = 10000
$ws.Range("K11").Value = 31800214
$ws.Range("L11").Value = 10000
$ws.Range("M11").Value = -31800074
$ws.Range("N11").Value = -10280
$ws.Range("H107").Value = 792.9032
$ws.Range("I107").Value = 659.2174
$ws.Range("K107").Value = 659.2174
$ws.Range("M107").Value = 1260.7826

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 499.15
$ws.Range("I107").Value = 317.25
$ws.Range("K107").Value = 317.25
$ws.Range("M107").Value = 1602.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1645.1
$ws.Range("I92").Value = 2249.5
$ws.Range("J92").Value = 738.5
$ws.Range("K92").Value = 6748.5
$ws.Range("L92").Value = 2215.5
$ws.Range("M92").Value = -5500.5
$ws.Range("N92").Value = -4711.5
$ws.Range("H107").Value = 1539.9857
$ws.Range("J107").Value = 2275.342
$ws.Range("L107").Value = 6826.026
$ws.Range("N107").Value = -10666.026

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 893.3570999999999
$ws.Range("I107").Value = 993.5454999999999
$ws.Range("K107").Value = 993.5454999999999
$ws.Range("M107").Value = 926.4545000000001
$ws.Range("H113").Value = 1768.1428
$ws.Range("I113").Value = 1479.5
$ws.Range("K113").Value = 1479.5
$ws.Range("M113").Value = 690.5
$ws.Range("H122").Value = 29413140
$ws.Range("I122").Value = 1452.0769
$ws.Range("J122").Value = 125001130
$ws.Range("K122").Value = 4356.2307
$ws.Range("L122").Value = 375003390
$ws.Range("M122").Value = -1906.2307
$ws.Range("N122").Value = -375008290
$ws.Range("H132").Value = 2817.111
$ws.Range("I132").Value = 2544.3125
$ws.Range("J132").Value = 4999.5
$ws.Range("K132").Value = 7632.9375
$ws.Range("L132").Value = 14998.5
$ws.Range("M132").Value = -5102.9375
$ws.Range("N132").Value = -20058.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1935.4
$ws.Range("I61").Value = 1935.4
$ws.Range("K61").Value = 1935.4
$ws.Range("M61").Value = -1733.4
$ws.Range("H68").Value = 1999
$ws.Range("I68").Value = 1999
$ws.Range("K68").Value = 1999
$ws.Range("M68").Value = -1250
$ws.Range("H71").Value = 1999
$ws.Range("I71").Value = 1999
$ws.Range("K71").Value = 9995
$ws.Range("M71").Value = -6251
$ws.Range("H100").Value = 3500.75
$ws.Range("I100").Value = 2001.5
$ws.Range("K100").Value = 2001.5
$ws.Range("M100").Value = -1460.5
$ws.Range("H110").Value = 52832.668
$ws.Range("J110").Value = 52832.668
$ws.Range("L110").Value = 52832.668
$ws.Range("N110").Value = -61012.668
$ws.Range("H113").Value = 1935.4
$ws.Range("I113").Value = 1935.4
$ws.Range("K113").Value = 1935.4
$ws.Range("M113").Value = 234.5999999999999
$ws.Range("H122").Value = 4240
$ws.Range("I122").Value = 3995
$ws.Range("K122").Value = 11985
$ws.Range("M122").Value = -9535
$ws.Range("H132").Value = 3712.923
$ws.Range("I132").Value = 3125.5
$ws.Range("J132").Value = 4513.9546
$ws.Range("K132").Value = 9376.5
$ws.Range("L132").Value = 13541.8638
$ws.Range("M132").Value = -6846.5
$ws.Range("N132").Value = -18601.8638
$ws.Range("H136").Value = 5598.2
$ws.Range("I136").Value = 2999.75
$ws.Range("K136").Value = 8999.25
$ws.Range("M136").Value = -6449.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 142860210
$ws.Range("I96").Value = 4000
$ws.Range("K96").Value = 4000
$ws.Range("M96").Value = -2627
$ws.Range("H132").Value = 27280.092
$ws.Range("I132").Value = 29728.166
$ws.Range("J132").Value = 2799.3333
$ws.Range("K132").Value = 89184.49800000001
$ws.Range("L132").Value = 8397.999899999999
$ws.Range("M132").Value = -86654.49800000001
$ws.Range("N132").Value = -13457.9999
$ws.Range("H136").Value = 9410.275
$ws.Range("I136").Value = 13255.272
$ws.Range("K136").Value = 39765.81600000001
$ws.Range("M136").Value = -37215.81600000001
